$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New injury-data rows (78-97) appended below the existing table.
$data = @(
  @{Row=78; A=38485; B="D229160"; C=$null;  D=$null; E=$null; F=$null; G=$null; H=$null},
  @{Row=79; A=38485; B="D229161"; C=$null;  D=$null; E=$null; F=$null; G=$null; H=$null},
  @{Row=80; A=38843; B="D229103"; C=1190;  D=5;   E=0;   F=0;   G=0.6; H=0},
  @{Row=81; A=38844; B="D229097"; C=1230;  D=10;  E=5;   F=40;  G=0;   H=0},
  @{Row=82; A=38847; B="D217910"; C=1290;  D=15;  E=10;  F=0;   G=1;   H=0},
  @{Row=83; A=39206; B="D229167"; C=1220;  D=15;  E=35;  F=50;  G=1.5; H=6},
  @{Row=84; A=39207; B="D229160"; C=1250;  D=15;  E=15;  F=0;   G=3.6; H=12},
  @{Row=85; A=39207; B="D229170"; C=1260;  D=20;  E=15;  F=5;   G=2.5; H=15},
  @{Row=86; A=39208; B="D229165"; C=1350;  D=30;  E=10;  F=5;   G=0.8; H=0},
  @{Row=87; A=39208; B="D229679"; C=1310;  D=15;  E=10;  F=0;   G=0;   H=3},
  @{Row=88; A=39208; B="D229735"; C=1330;  D=15;  E=10;  F=5;   G=0;   H=0},
  @{Row=89; A=39209; B="D229658"; C=1180;  D=0;   E=0;   F=0;   G=0.3; H=0},
  @{Row=90; A=39209; B="D229025"; C=1170;  D=10;  E=5;   F=5;   G=0;   H=0},
  @{Row=91; A=39211; B="D229776"; C=1190;  D=25;  E=20;  F=0;   G=8;   H=2.4},
  @{Row=92; A=39211; B="D237389"; C=1200;  D=10;  E=20;  F=0;   G=0;   H=0},
  @{Row=93; A=39212; B="D229775"; C=1200;  D=30;  E=30;  F=0;   G=0;   H=2.4},
  @{Row=94; A=39217; B="D237278"; C=1210;  D=15;  E=20;  F=60;  G=0.9; H=4.2},
  @{Row=95; A=39217; B="D237361"; C=1100;  D=5;   E=0;   F=0;   G=0;   H=3.6},
  @{Row=96; A=39217; B="D237378"; C=1160;  D=0;   E=5;   F=0;   G=1.2; H=4.8},
  @{Row=97; A=39217; B="D229162"; C=1220;  D=25;  E=15;  F=0;   G=7;   H=2.4}
)

foreach ($r in $data) {
  $ws.Cells.Item($r.Row, 1).Value2 = $r.A
  $ws.Cells.Item($r.Row, 2).Value = $r.B
  if ($r.C -ne $null) { $ws.Cells.Item($r.Row, 3).Value2 = $r.C }
  if ($r.D -ne $null) { $ws.Cells.Item($r.Row, 4).Value2 = $r.D }
  if ($r.E -ne $null) { $ws.Cells.Item($r.Row, 5).Value2 = $r.E }
  if ($r.F -ne $null) { $ws.Cells.Item($r.Row, 6).Value2 = $r.F }
  if ($r.G -ne $null) { $ws.Cells.Item($r.Row, 7).Value2 = $r.G }
  if ($r.H -ne $null) { $ws.Cells.Item($r.Row, 8).Value2 = $r.H }
}

# Give the new dates (col A) the same plain short-date format used by the
# rest of the sheet's date column (no fill), then stamp it across the
# whole new range from a single source cell so they all share one style.
$ws.Range("A78").NumberFormat = "mm-dd-yy"
$ws.Range("A78").Copy()
$ws.Range("A79:A97").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the sheet view down to the newly added rows.
$ws.Application.ActiveWindow.ScrollRow = 80
$ws.Range("A82").Select()
